$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep Price/Volume (and Coin/Link) cells as plain text so Excel does not
# reinterpret numeric-looking strings (e.g. "214.77", "1.00") as numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.502.89"
$ws.Range("E2").Value = "  +4.76%  "
$ws.Range("D3").Value = "1.590.79"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").Value = "214.77"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").Value = "0.498"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("D8").Value = "23.89"
$ws.Range("E8").Value = "  +8.49%  "
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").Value = "0.0601"
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").Value = "1.818.02"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").Value = "1.608.35"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "0.531"
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("D16").Value = "28.464.94"
$ws.Range("E16").Value = "  +4.82%  "
$ws.Range("D17").Value = "64.04"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("D18").Value = "232.58"
$ws.Range("E18").Value = "  +7.31%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "0.0₃0710"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "4.14"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "9.42"
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D25").Value = "151.90"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "0.108"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "1.419.39"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("E36").Value = "  -6.13%  "
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").Value = "2.56"
$ws.Range("E39").Value = "  +9.45%  "
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("D41").Value = "0.813"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.76"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "1.84"
$ws.Range("E44").Value = "  +6.62%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "0.979"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").Value = "64.63"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "1.728.51"
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("D48").Value = "87.90"
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("D49").Value = "0.0⁦0106"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("D50").Value = "0.0523"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "39.59"
$ws.Range("E51").Value = "  +16.71%  "
